# Update the sample "master-machine_spec" rows from the SanDisk Cruzer Blade
# USB-drive sample data to the Dell Vostro computer sample data (with its
# Arabic translations), as per the "Adding Master Data XLS" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (id 1001, English / eng record) ---------------------------------
$ws.Range("B2").Value = "Vostro"
$ws.Range("C2").Value = "Dell"
$ws.Range("D2").Value = 3568
$ws.Range("E2").Value = "DKS"
$ws.Range("F2").Value = 1.454
$ws.Range("G2").Value = "To take enrollments"
$ws.Range("H2").Value = "eng"
$ws.Range("J2").Value = "superadmin"
$ws.Range("K2").Value = "now()"

# --- Row 3 (id 1002, Arabic / ara record) -----------------------------------
$ws.Range("B3").Value = "ستر  "
$ws.Range("C3").Value = "دلّ  "
$ws.Range("D3").Value = 3568
$ws.Range("E3").Value = "DKS"
$ws.Range("F3").Value = 1.454
$ws.Range("G3").Value = "لأخذ التسجيلات"
$ws.Range("H3").Value = "ara"
$ws.Range("J3").Value = "superadmin"
$ws.Range("K3").Value = "now()"

# --- View state: select whole rows 4 downwards, as in the saved file -------
$ws.Rows("4:1048576").Select() | Out-Null

# --- Page setup used for printing the sheet ---------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
